$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1001.1818
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 1071.3
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 3213.9
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -3549.9

$ws.Range("H113").Value = 2141
$ws.Range("I113").Value = 2174.75
$ws.Range("K113").Value = 2174.75
$ws.Range("M113").Value = 1079.25

$ws.Range("H137").Value = 1382.8445
$ws.Range("I137").Value = 1015.44446
$ws.Range("J137").Value = 1933.9445
$ws.Range("K137").Value = 3046.33338
$ws.Range("L137").Value = 5801.833500000001
$ws.Range("M137").Value = -496.33338
$ws.Range("N137").Value = -10901.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 2504000
$ws.Range("J92").Value = 2504000
$ws.Range("L92").Value = 2504000
$ws.Range("N92").Value = -2508992

$ws.Range("H122").Value = 2677.5
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 2345.7144
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 7037.1432
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -11937.1432

$ws.Range("H132").Value = 2815.9333
$ws.Range("I132").Value = 2411.4546
$ws.Range("K132").Value = 7234.3638
$ws.Range("M132").Value = -4704.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 333334660
$ws.Range("I105").Value = 333334660
$ws.Range("K105").Value = 333334660
$ws.Range("M105").Value = -333332913

$ws.Range("H134").Value = 2832.585
$ws.Range("I134").Value = 652.4211
$ws.Range("J134").Value = 8355.666999999999
$ws.Range("K134").Value = 1957.2633
$ws.Range("L134").Value = 25067.001
$ws.Range("M134").Value = 577.7366999999999
$ws.Range("N134").Value = -30137.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 781.1539
$ws.Range("I122").Value = 841.36365
$ws.Range("K122").Value = 2524.09095
$ws.Range("M122").Value = -74.09094999999979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 135.8
$ws.Range("J2").Value = 190
$ws.Range("L2").Value = 1140
$ws.Range("N2").Value = -1366

$ws.Range("H12").Value = 141.94444
$ws.Range("I12").Value = 349.66666
$ws.Range("J12").Value = 100.4
$ws.Range("K12").Value = 1048.99998
$ws.Range("L12").Value = 301.2
$ws.Range("M12").Value = -875.9999800000001
$ws.Range("N12").Value = -647.2

$ws.Range("H24").Value = 1002
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1002
$ws.Range("K24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("M24").Value = 3006
$ws.Range("N24").Value = -3466

$ws.Range("H74").Value = 5600
$ws.Range("J74").Value = 5600
$ws.Range("L74").Value = 16800
$ws.Range("N74").Value = -18922

$ws.Range("H77").Value = 5600
$ws.Range("J77").Value = 5600
$ws.Range("L77").Value = 50400
$ws.Range("N77").Value = -61008

$ws.Range("H87").Value = 1823.7778
$ws.Range("J87").Value = 2000
$ws.Range("L87").Value = 6000
$ws.Range("N87").Value = -8496

$ws.Range("H88").Value = 7045.4546
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H90").Value = 1823.7778
$ws.Range("J90").Value = 2000
$ws.Range("L90").Value = 18000
$ws.Range("N90").Value = -30480

$ws.Range("H91").Value = 7045.4546
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H126").Value = 2686.6
$ws.Range("I126").Value = 1850
$ws.Range("K126").Value = 5550
$ws.Range("M126").Value = -610

$ws.Range("H129").Value = 32053190
$ws.Range("I129").Value = 111111730
$ws.Range("J129").Value = 8335630.5
$ws.Range("K129").Value = 333335190
$ws.Range("L129").Value = 25006891.5
$ws.Range("M129").Value = -333330190
$ws.Range("N129").Value = -25016891.5

$ws.Range("H130").Value = 2315.9167
$ws.Range("I130").Value = 1030
$ws.Range("J130").Value = 2432.818
$ws.Range("K130").Value = 3090
$ws.Range("L130").Value = 7298.454000000001
$ws.Range("M130").Value = 1930
$ws.Range("N130").Value = -17338.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5142.1665
$ws.Range("I80").Value = 4763.125
$ws.Range("J80").Value = 5900.25
$ws.Range("K80").Value = 4763.125
$ws.Range("L80").Value = 5900.25
$ws.Range("M80").Value = -3765.125
$ws.Range("N80").Value = -7896.25

$ws.Range("H83").Value = 5142.1665
$ws.Range("I83").Value = 4763.125
$ws.Range("J83").Value = 5900.25
$ws.Range("K83").Value = 23815.625
$ws.Range("L83").Value = 29501.25
$ws.Range("M83").Value = -18823.625
$ws.Range("N83").Value = -39485.25

$ws.Range("H97").Value = 730
$ws.Range("I97").Value = 730
$ws.Range("K97").Value = 730
$ws.Range("M97").Value = -234

$ws.Range("H122").Value = 8930408
$ws.Range("I122").Value = 1906.1666
$ws.Range("K122").Value = 5718.4998
$ws.Range("M122").Value = -3268.4998

$ws.Range("H126").Value = 2279.3333
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2878.5
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 8635.5
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -13575.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2101.25
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H40").Value = 2384.6843
$ws.Range("I40").Value = 2262.75
$ws.Range("J40").Value = 3035
$ws.Range("K40").Value = 2262.75
$ws.Range("L40").Value = 3035
$ws.Range("M40").Value = -2126.75
$ws.Range("N40").Value = -3307

$ws.Range("H126").Value = 2101.25
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws.Range("H132").Value = 2749
$ws.Range("I132").Value = 2854
$ws.Range("J132").Value = 2681.0588
$ws.Range("K132").Value = 8562
$ws.Range("L132").Value = 8043.176399999999
$ws.Range("M132").Value = -6032
$ws.Range("N132").Value = -13103.1764

$ws.Range("H136").Value = 1802.1666
$ws.Range("I136").Value = 1652
$ws.Range("J136").Value = 2102.5
$ws.Range("K136").Value = 4956
$ws.Range("L136").Value = 6307.5
$ws.Range("M136").Value = -2406
$ws.Range("N136").Value = -11407.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12501104
$ws.Range("I122").Value = 14707075
$ws.Range("K122").Value = 44121225
$ws.Range("M122").Value = -44118775

$ws.Range("H126").Value = 83334190
$ws.Range("I126").Value = 90909850
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 272729550
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -272727080
